# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2-132) from 2023-09-19 (serial 45188) to 2023-09-20 (serial 45189).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C132").Value = 45189
